# Append new rows (19-30) of build/test-run data to Sheet1, mirroring the
# "Updated json with build number" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be stored as text even when the literal content
    # looks like a number or a date (e.g. "81", "2021-03-31").
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Rows 19-24: six more "DAM Smoke Tests UAT" FAILURE runs ---------------
foreach ($r in 19..24) {
    $ws.Range("A$r").Value = "DAM Smoke Tests UAT"
    $ws.Range("D$r").Value = "01:05:17.790"
    $ws.Range("E$r").Value = "FAILURE"
    $ws.Range("F$r").Value = 72
    $ws.Range("G$r").Value = 71
    $ws.Range("H$r").Value = 6
    $ws.Range("I$r").Value = 65
    $ws.Range("J$r").Value = 0
    $ws.Range("K$r").Value = 0
}

# --- Row 25: "DAM Custom Execution" SUCCESS run -----------------------------
$ws.Range("A25").Value = "DAM Custom Execution"
$ws.Range("D25").Value = "00:01:20.396"
$ws.Range("E25").Value = "SUCCESS"
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 5
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0

# --- Row 26: "DAM Custom Execution" SUCCESS run, Build Number stored as text
$ws.Range("A26").Value = "DAM Custom Execution"
Set-TextValue $ws.Range("B26") "81"
Set-TextValue $ws.Range("C26") "2021-03-31"
$ws.Range("D26").Value = "00:01:32.984"
$ws.Range("E26").Value = "SUCCESS"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 10
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0

# --- Row 27: "DAM Custom Execution" SUCCESS run, Build Number numeric ------
$ws.Range("A27").Value = "DAM Custom Execution"
$ws.Range("B27").Value = 82
Set-TextValue $ws.Range("C27") "2021-03-31"
$ws.Range("D27").Value = "00:01:13.715"
$ws.Range("E27").Value = "SUCCESS"
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0

# --- Row 28: "DAM Smoke Tests UAT" FAILURE run, Build Number numeric -------
$ws.Range("A28").Value = "DAM Smoke Tests UAT"
$ws.Range("B28").Value = 13
Set-TextValue $ws.Range("C28") "2021-02-25"
$ws.Range("D28").Value = "01:05:17.790"
$ws.Range("E28").Value = "FAILURE"
$ws.Range("F28").Value = 72
$ws.Range("G28").Value = 71
$ws.Range("H28").Value = 6
$ws.Range("I28").Value = 65
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0

# --- Row 29: "DAM Custom Execution" SUCCESS run, Build Number numeric ------
$ws.Range("A29").Value = "DAM Custom Execution"
$ws.Range("B29").Value = 83
Set-TextValue $ws.Range("C29") "2021-03-31"
$ws.Range("D29").Value = "00:01:07.364"
$ws.Range("E29").Value = "SUCCESS"
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0

# --- Row 30: "DAM Custom Execution" SUCCESS run, Build Number numeric ------
$ws.Range("A30").Value = "DAM Custom Execution"
$ws.Range("B30").Value = 84
Set-TextValue $ws.Range("C30") "2021-03-31"
$ws.Range("D30").Value = "00:01:07.256"
$ws.Range("E30").Value = "SUCCESS"
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
